# Case_6_40 diagnostic.xlsx edit
# Adds a small 2x2 "disconnected_elements" diagnostic block to Sheet1:
#   B1 = 0            (bold, centered/top, thin box border)
#   A2 = 0            (bold, centered/top, thin box border)
#   B2 = "disconnected_elements" (shared string, no special formatting)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- Shared formatting (bold font, thin box border, center/top align) -
# Build the format once on a scratch cell, then copy/paste-special the
# format onto both target cells so they end up sharing a single cell
# style (xf) entry instead of each accumulating its own intermediate
# style snapshot.
$template = $ws.Range("Z1")
$template.Font.Bold = $true
$template.HorizontalAlignment = -4108
$template.VerticalAlignment = -4160
$template.Borders.LineStyle = 1
$template.Borders.Weight = 2

$template.Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("A2").PasteSpecial(-4122)

$template.Clear()
